# Align figures with ORBIT and reformat the workbook's active selections.
$wb = $excel.ActiveWorkbook

$wsPort = $wb.Worksheets.Item("Port-investments")
$wsShort = $wb.Worksheets.Item("schedule-short")
$wsSched = $wb.Worksheets.Item("schedule")

# --- Port-investments: updated investment figures (row 7-11) ---
# Row 7
$wsPort.Range("K7").Value = 700
$wsPort.Range("L7").Value = 200
$wsPort.Range("M7").Value = 350
$wsPort.Range("N7").Value = 600

# Row 8
$wsPort.Range("L8").Value = 125
$wsPort.Range("M8").Value = 250
$wsPort.Range("N8").Value = 500

# Row 9
$wsPort.Range("L9").Value = 70
$wsPort.Range("M9").Value = 200
$wsPort.Range("N9").Value = 400

# Row 10
$wsPort.Range("K10").Value = 1100
$wsPort.Range("N10").Value = 35

# Row 11
$wsPort.Range("K11").Value = 665
$wsPort.Range("L11").Value = 0
$wsPort.Range("M11").Value = 0
$wsPort.Range("N11").Value = 0

# --- View / selection state per sheet ---
$wsPort.Activate()
$wsPort.Range("L21").Select()

$wsShort.Activate()
$wsShort.Range("G10").Select()

# schedule becomes the active sheet/tab on save
$wsSched.Activate()
$wsSched.Range("E38").Select()
